# Insert a new daily price record as the new row 38 ("Hortaliza, Vega Modelo
# de Temuco - Espinaca"). This pushes the existing rows 38:119 down to
# 39:120 (dimension grows from A1:R119 to A1:R120), matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 38, shifting rows 38-119 down to 39-120.
$ws.Rows("38:38").Insert()

# Populate the newly inserted row 38 with the new record's data.
$ws.Range("A38").Value = 10
$ws.Range("B38").Value = "Vega Modelo de Temuco"
$ws.Range("C38").Value = "La Araucanía"
$ws.Range("D38").Value = 44581
$ws.Range("E38").Value = 9
$ws.Range("F38").Value = 100112012
$ws.Range("G38").Value = "Espinaca"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 70
$ws.Range("K38").Value = 13000
$ws.Range("L38").Value = 14000
$ws.Range("M38").Value = 13571
$ws.Range("N38").Value = "`$/docena de atados"
$ws.Range("O38").Value = "Región de La Araucanía"
$ws.Range("P38").Value = 4524
$ws.Range("Q38").Value = 3
$ws.Range("R38").Value = "Hortaliza"
